# Update Detailed run sheet optimization
# - Constraints sheet: Clinic_Close time 20:00 -> 23:00, widen the Value column
# - Worker sheet: reorganise Oral Therapist / Dentist job codes and add new staff

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Constraints sheet: Clinic_Close (row 14) value changes from 20:00 to 23:00
# ---------------------------------------------------------------------------
$wsConstraints = $wb.Worksheets.Item("Constraints")
$wsConstraints.Cells.Item(14, 2).Value = 0.95833333333333337
$wsConstraints.Columns.Item(2).ColumnWidth = 10.9

# ---------------------------------------------------------------------------
# Worker sheet: rebuild the staff roster from row 4 onward
# ---------------------------------------------------------------------------
$wsWorker = $wb.Worksheets.Item("Worker")

$workers = @(
    @("Klarion",  "OHT1", "OHT", 45000),
    @("Cindy",    "OHT2", "OHT", 45000),
    @("Arlong",   "OHT3", "OHT", 48000),
    @("Marlo",    "OHT4", "OHT", 48000),
    @("Reno",     "OHT5", "OHT", 49000),
    @("Lorean",   "OHT6", "OHT", 49000),
    @("Jannik",   "OHT7", "OHT", 50000),
    @("Shane",    "D1N1", "LV1", 56000),
    @("Mary",     "D1N2", "LV1", 62000),
    @("Nashvile", "D1N3", "LV1", 65000),
    @("Larry",    "D1N4", "LV1", 65000),
    @("Meghan",   "D2N1", "LV2", 68000),
    @("Sharon",   "D2N2", "LV2", 70000),
    @("Ello",     "D2N3", "LV2", 75000),
    @("Myrion",   "D3N1", "LV3", 80000),
    @("Arthur",   "D3N2", "LV3", 85000),
    @("Darius",   "D3N3", "LV3", 90000)
)

$r = 4
foreach ($row in $workers) {
    $wsWorker.Cells.Item($r, 1).Value = $row[0]
    $wsWorker.Cells.Item($r, 2).Value = $row[1]
    $wsWorker.Cells.Item($r, 3).Value = $row[2]
    $wsWorker.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$wsWorker.Range("D21").Select() | Out-Null

Write-Output "Worker roster updated; Constraints updated"
